$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column D ("Price") holds numeric-looking text (e.g. "241.79", "0.9130").
# A plain Range.Value assignment lets Excel auto-convert those into real
# numbers (dropping formatting like trailing zeros), so we briefly force a
# Text format while writing, then restore the default "Normal" style so the
# cell keeps matching its original (unstyled) look.
function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '29.508.62'
$ws.Range('E2').Value = '  +1.04%  '
Set-TextValue 'D3' '1.880.90'
$ws.Range('E3').Value = '  +1.55%  '
$ws.Range('E4').Value = '  +0.01%  '
Set-TextValue 'D5' '0.7152'
$ws.Range('E5').Value = '  +2.64%  '
Set-TextValue 'D6' '241.79'
$ws.Range('E7').Value = '  +0.06%  '
Set-TextValue 'D8' '0.07937'
$ws.Range('E8').Value = '  +1.22%  '
Set-TextValue 'D9' '0.3111'
$ws.Range('E9').Value = '  +3.16%  '
Set-TextValue 'D10' '25.26'
$ws.Range('E10').Value = '  +6.60%  '
Set-TextValue 'D11' '0.08275'
$ws.Range('E11').Value = '  +2.02%  '
Set-TextValue 'D12' '0.7296'
$ws.Range('E12').Value = '  +3.72%  '
Set-TextValue 'D13' '5.286'
$ws.Range('E13').Value = '  +2.24%  '
Set-TextValue 'D14' '1.852.99'
$ws.Range('E14').Value = '  +0.02%  '
Set-TextValue 'D15' '91.23'
$ws.Range('E15').Value = '  +2.01%  '
Set-TextValue 'D16' '29.527.12'
$ws.Range('E16').Value = '  +1.06%  '
Set-TextValue 'D17' '5.937'
$ws.Range('E17').Value = '  +2.39%  '
Set-TextValue 'D18' '246.50'
$ws.Range('E18').Value = '  +4.64%  '
Set-TextValue 'D19' '0.000007895'
$ws.Range('E19').Value = '  +1.06%  '
$ws.Range('E20').Value = '  +1.35%  '
Set-TextValue 'D21' '2.129.63'
$ws.Range('E21').Value = '  +1.39%  '
$ws.Range('E22').Value = '  +0.07%  '
Set-TextValue 'D23' '7.975'
$ws.Range('E23').Value = '  +6.30%  '
$ws.Range('E24').Value = '  -0.02%  '
Set-TextValue 'D25' '0.1618'
$ws.Range('E25').Value = '  +14.34%  '
Set-TextValue 'D26' '163.19'
$ws.Range('E26').Value = '  +0.42%  '
Set-TextValue 'D27' '9.071'
$ws.Range('E27').Value = '  +2.41%  '
Set-TextValue 'D28' '18.33'
$ws.Range('E28').Value = '  +1.74%  '
Set-TextValue 'D29' '1.357'
$ws.Range('E29').Value = '  -3.34%  '
Set-TextValue 'D30' '1.495'
$ws.Range('E30').Value = '  +1.18%  '
Set-TextValue 'D31' '4.388'
$ws.Range('E31').Value = '  +1.63%  '
Set-TextValue 'D32' '4.109'
$ws.Range('E32').Value = '  +2.56%  '
Set-TextValue 'D33' '0.05278'
$ws.Range('E33').Value = '  +2.57%  '
Set-TextValue 'D34' '1.964'
$ws.Range('E34').Value = '  +2.25%  '
$ws.Range('E35').Value = '  +3.17%  '
Set-TextValue 'D36' '0.7275'
$ws.Range('E36').Value = '  +2.62%  '
Set-TextValue 'D37' '2.678'
$ws.Range('E37').Value = '  -0.01%  '
Set-TextValue 'D38' '0.01869'
Set-TextValue 'D39' '1.231.52'
$ws.Range('E39').Value = '  +6.98%  '
Set-TextValue 'D40' '2.716'
$ws.Range('E40').Value = '  +0.45%  '
Set-TextValue 'D41' '0.9130'
Set-TextValue 'D42' '6.196'
$ws.Range('E42').Value = '  +3.90%  '
$ws.Range('E43').Value = '  +5.70%  '
$ws.Range('E44').Value = '  +0.08%  '
Set-TextValue 'D45' '102.35'
$ws.Range('E45').Value = '  -0.57%  '
Set-TextValue 'D46' '2.024.51'
$ws.Range('E46').Value = '  +1.36%  '
Set-TextValue 'D47' '0.5287'
$ws.Range('E47').Value = '  -0.20%  '
Set-TextValue 'D48' '1.805'
$ws.Range('E48').Value = '  +4.16%  '
Set-TextValue 'D49' '2.939'
$ws.Range('E49').Value = '  +11.29%  '
Set-TextValue 'D50' '0.00000000121'
$ws.Range('E50').Value = '  +1.80%  '
Set-TextValue 'D51' '9.315'
$ws.Range('E51').Value = '  +2.03%  '
